$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

$cells = @("B11","C11","B12","C12","D12","B14","C14","D14")
foreach ($addr in $cells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Enterprises density (per 1000 people)
$ws.Range("B11").Value = "24.89"
$ws.Range("C11").Value = "0.92"

# Employment (% of total)
$ws.Range("B12").Value = "50.81"
$ws.Range("C12").Value = "23.21"
$ws.Range("D12").Value = "74.02"

# Enterprises (% of total)
$ws.Range("B14").Value = "96.16"
$ws.Range("C14").Value = "3.54"
$ws.Range("D14").Value = "99.69"
